# Auto-generated script to update market-price snapshot values
# across multiple Leve-profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 9842.333000000001
$ws.Range("I69").Value = 1999.5
$ws.Range("J69").Value = 13763.75
$ws.Range("K69").Value = 5998.5
$ws.Range("L69").Value = 41291.25
$ws.Range("M69").Value = -5124.5
$ws.Range("N69").Value = -43039.25
$ws.Range("H70").Value = 3009.4
$ws.Range("J70").Value = 3021.7778
$ws.Range("L70").Value = 9065.3334
$ws.Range("N70").Value = -9605.3334
$ws.Range("H72").Value = 9842.333000000001
$ws.Range("I72").Value = 1999.5
$ws.Range("J72").Value = 13763.75
$ws.Range("K72").Value = 17995.5
$ws.Range("L72").Value = 123873.75
$ws.Range("M72").Value = -13627.5
$ws.Range("N72").Value = -132609.75
$ws.Range("H73").Value = 3009.4
$ws.Range("J73").Value = 3021.7778
$ws.Range("L73").Value = 9065.3334
$ws.Range("N73").Value = -10937.3334
$ws.Range("H80").Value = 931.7308
$ws.Range("I80").Value = 687.5833
$ws.Range("J80").Value = 1141
$ws.Range("K80").Value = 2062.7499
$ws.Range("L80").Value = 3423
$ws.Range("M80").Value = -1064.7499
$ws.Range("N80").Value = -5419
$ws.Range("H83").Value = 931.7308
$ws.Range("I83").Value = 687.5833
$ws.Range("J83").Value = 1141
$ws.Range("K83").Value = 6188.2497
$ws.Range("L83").Value = 10269
$ws.Range("M83").Value = -1196.2497
$ws.Range("N83").Value = -20253
$ws.Range("H92").Value = 568.75
$ws.Range("I92").Value = 568.75
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 568.75
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = 679.25
$ws.Range("H107").Value = 2114.7856
$ws.Range("I107").Value = 1095
$ws.Range("J107").Value = 2879.625
$ws.Range("K107").Value = 1095
$ws.Range("L107").Value = 2879.625
$ws.Range("M107").Value = 825
$ws.Range("N107").Value = -6719.625
$ws.Range("H137").Value = 7486.607
$ws.Range("I137").Value = 7616.081
$ws.Range("K137").Value = 22848.243
$ws.Range("M137").Value = -20298.243
$ws.Range("H138").Value = 4746.9287
$ws.Range("I138").Value = 2722.375
$ws.Range("J138").Value = 5556.75
$ws.Range("K138").Value = 8167.125
$ws.Range("L138").Value = 16670.25
$ws.Range("M138").Value = -3027.125
$ws.Range("N138").Value = -26950.25

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1688.8718
$ws.Range("I32").Value = 1449.4722
$ws.Range("K32").Value = 1449.4722
$ws.Range("M32").Value = -1162.4722
$ws.Range("H46").Value = 1250
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1250
$ws.Range("K46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("M46").Value = 1250
$ws.Range("N46").Value = -1888
$ws.Range("H74").Value = 80257.57000000001
$ws.Range("I74").Value = 101271.4
$ws.Range("J74").Value = 27723
$ws.Range("K74").Value = 101271.4
$ws.Range("L74").Value = 27723
$ws.Range("M74").Value = -100397.4
$ws.Range("N74").Value = -29471
$ws.Range("H77").Value = 80257.57000000001
$ws.Range("I77").Value = 101271.4
$ws.Range("J77").Value = 27723
$ws.Range("K77").Value = 506357
$ws.Range("L77").Value = 138615
$ws.Range("M77").Value = -501989
$ws.Range("N77").Value = -147351

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1386.4814
$ws.Range("I94").Value = 1346.826
$ws.Range("J94").Value = 1614.5
$ws.Range("K94").Value = 1346.826
$ws.Range("L94").Value = 1614.5
$ws.Range("M94").Value = -895.826
$ws.Range("N94").Value = -2516.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 189.44444
$ws.Range("I5").Value = 86
$ws.Range("J5").Value = 318.75
$ws.Range("K5").Value = 86
$ws.Range("L5").Value = 318.75
$ws.Range("M5").Value = 26
$ws.Range("N5").Value = -542.75
$ws.Range("H15").Value = 1601.8572
$ws.Range("I15").Value = 141
$ws.Range("K15").Value = 141
$ws.Range("M15").Value = 29
$ws.Range("H26").Value = 8182
$ws.Range("J26").Value = 8182
$ws.Range("L26").Value = 8182
$ws.Range("N26").Value = -8756
$ws.Range("H29").Value = 14775.6
$ws.Range("J29").Value = 14775.6
$ws.Range("L29").Value = 14775.6
$ws.Range("N29").Value = -15361.6
$ws.Range("H31").Value = 4874901
$ws.Range("I31").Value = 6665171
$ws.Range("K31").Value = 6665171
$ws.Range("M31").Value = -6664876
$ws.Range("H33").Value = 34426.285
$ws.Range("J33").Value = 45596.8
$ws.Range("L33").Value = 45596.8
$ws.Range("N33").Value = -46354.8
$ws.Range("H34").Value = 4874901
$ws.Range("I34").Value = 6665171
$ws.Range("K34").Value = 6665171
$ws.Range("M34").Value = -6664969
$ws.Range("H38").Value = 3979.3333
$ws.Range("I38").Value = 4269
$ws.Range("K38").Value = 4269
$ws.Range("M38").Value = -3892
$ws.Range("H46").Value = 3979.3333
$ws.Range("I46").Value = 4269
$ws.Range("K46").Value = 4269
$ws.Range("M46").Value = -4058

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 65.15385000000001
$ws.Range("I38").Value = 73.933334
$ws.Range("J38").Value = 53.18182
$ws.Range("K38").Value = 221.800002
$ws.Range("L38").Value = 159.54546
$ws.Range("M38").Value = 125.199998
$ws.Range("N38").Value = -853.54546
$ws.Range("H68").Value = 1790039.8
$ws.Range("I68").Value = 4356.6
$ws.Range("J68").Value = 2178231.8
$ws.Range("K68").Value = 13069.8
$ws.Range("L68").Value = 6534695.399999999
$ws.Range("M68").Value = -12258.8
$ws.Range("N68").Value = -6536317.399999999
$ws.Range("H71").Value = 1790039.8
$ws.Range("I71").Value = 4356.6
$ws.Range("J71").Value = 2178231.8
$ws.Range("K71").Value = 39209.4
$ws.Range("L71").Value = 19604086.2
$ws.Range("M71").Value = -35153.4
$ws.Range("N71").Value = -19612198.2
$ws.Range("H129").Value = 19803364
$ws.Range("I129").Value = 24751698
$ws.Range("K129").Value = 74255094
$ws.Range("M129").Value = -74250094
$ws.Range("H131").Value = 22134.47
$ws.Range("I131").Value = 92216.09
$ws.Range("K131").Value = 276648.27
$ws.Range("M131").Value = -271608.27
$ws.Range("H137").Value = 6669.8887
$ws.Range("J137").Value = 4699.5
$ws.Range("L137").Value = 14098.5
$ws.Range("N137").Value = -24298.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 20498.053
$ws.Range("I126").Value = 27417.615
$ws.Range("J126").Value = 5505.6665
$ws.Range("K126").Value = 82252.845
$ws.Range("L126").Value = 16516.9995
$ws.Range("M126").Value = -79782.845
$ws.Range("N126").Value = -21456.9995

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 22140.916
$ws.Range("I41").Value = 17447
$ws.Range("J41").Value = 23705.555
$ws.Range("K41").Value = 17447
$ws.Range("L41").Value = 23705.555
$ws.Range("M41").Value = -17057
$ws.Range("N41").Value = -24485.555
$ws.Range("H62").Value = 20711.904
$ws.Range("I62").Value = 19997.37
$ws.Range("J62").Value = 27500
$ws.Range("K62").Value = 19997.37
$ws.Range("L62").Value = 27500
$ws.Range("M62").Value = -19373.37
$ws.Range("N62").Value = -28748
$ws.Range("H65").Value = 20711.904
$ws.Range("I65").Value = 19997.37
$ws.Range("J65").Value = 27500
$ws.Range("K65").Value = 99986.84999999999
$ws.Range("L65").Value = 137500
$ws.Range("M65").Value = -96866.84999999999
$ws.Range("N65").Value = -143740
$ws.Range("H81").Value = 10274.889
$ws.Range("I81").Value = 36316
$ws.Range("J81").Value = 5066.6665
$ws.Range("K81").Value = 72632
$ws.Range("L81").Value = 10133.333
$ws.Range("M81").Value = -71571
$ws.Range("N81").Value = -12255.333
$ws.Range("H84").Value = 10274.889
$ws.Range("I84").Value = 36316
$ws.Range("J84").Value = 5066.6665
$ws.Range("K84").Value = 363160
$ws.Range("L84").Value = 50666.665
$ws.Range("M84").Value = -357856
$ws.Range("N84").Value = -61274.665
